# Add new rows 227-234 (columns C:E) to Sheet1, mirroring the pattern of
# the existing table, and update the sheet view to the final scroll/selection
# state (topLeftCell A214, active cell G225).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style used for F1 (and, per the target, for most of the new cells) -
# grab it once so we can best-effort carry it onto the new cells below.
$f1Style = $ws.Range("F1").Style

# Row 227 - plain (no explicit style)
$ws.Range("C227").Value = 26.15
$ws.Range("D227").Value = 133
$ws.Range("E227").Value = 84

# Row 228 - only C228 carries the F1-style
$ws.Range("C228").Value = 26.15
$ws.Range("C228").Style = $f1Style
$ws.Range("D228").Value = 132
$ws.Range("E228").Value = 83

# Rows 229-234 - alternating 133/84 and 132/83, all three cells styled
$data = @(
    @(229, 133, 84),
    @(230, 132, 83),
    @(231, 133, 84),
    @(232, 132, 83),
    @(233, 133, 84),
    @(234, 132, 83)
)

foreach ($row in $data) {
    $r = $row[0]
    $dVal = $row[1]
    $eVal = $row[2]

    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Value = 26.15
    $cCell.Style = $f1Style

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = $dVal
    $dCell.Style = $f1Style

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value = $eVal
    $eCell.Style = $f1Style
}

# Final UI state: selection on G225 (also nudges topLeftCell scrolling
# toward the bottom of the newly extended table).
$ws.Range("G225").Select()
